$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$rng = $ws.Range("A1")
Write-Output $rng.Font.Name
Write-Output $rng.Font.Size
Write-Output $rng.Font.Bold
Write-Output $rng.HorizontalAlignment
Write-Output $rng.VerticalAlignment
Write-Output $rng.WrapText
Write-Output $rng.Interior.Color
Write-Output $rng.Borders.Item(7).LineStyle
